$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column L (12) -- existing L:P shift to M:Q
$ws.Columns.Item(12).Insert()

# Copy the header formatting (bold/fill style) from the neighboring header
# cell (K1) onto the new header cell, then set its text.
$ws.Cells.Item(1, 11).Copy()
$ws.Cells.Item(1, 12).PasteSpecial(-4122)
$ws.Cells.Item(1, 12).Value = "ORDEN"

# Scroll the view over and select the new column's data cell, like the
# author did while reviewing the freshly inserted "ORDEN" column.
$excel.Goto($ws.Range("L4"), $true)

$excel.CutCopyMode = $false
